# Automatic update of files.
#
# 1) Column C ("Förändrad") gets bumped from 45184 to 45186 for every data
#    row (rows 2-260).
# 2) The HYPERLINK() formulas in columns S,T,U,V,W,X,Y for rows 2-33 gain a
#    second "friendly name" argument equal to the row's "Beteckning"
#    (column A) value, e.g.
#      HYPERLINK("...A 20688-2022.xlsx")
#    becomes
#      HYPERLINK("...A 20688-2022.xlsx", "A 20688-2022")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 260
$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = $firstRow; $r -le $lastRow; $r++) {

    # --- 1) bump the "Förändrad" date in column C ---
    $ws.Range("C" + $r).Value = 45186

    # --- 2) add the friendly-name argument to every HYPERLINK formula on this row ---
    $beteckning = $ws.Range("A" + $r).Value2

    if (-not [string]::IsNullOrEmpty($beteckning)) {
        foreach ($col in $hyperlinkCols) {
            $cell = $ws.Range($col + $r)
            $formula = $cell.Formula

            if (-not [string]::IsNullOrEmpty($formula)) {
                # only touch plain single-argument HYPERLINK("...") formulas;
                # leave any already-updated ones untouched
                if ($formula -match '^=HYPERLINK\("[^"]*"\)$') {
                    $newFormula = $formula -replace '\)$', (', "' + $beteckning + '")')
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
